$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit applies a permutation of match-record rows: the content of
# several rows (columns B through AC) is rotated among rows, while the
# sequential index in column A stays fixed per row position.

# Snapshot the current (pre-edit) B:AC values for every row involved,
# BEFORE any writes happen, so the rotations do not clobber each other.
$snap = @{}
$snap[3] = $ws.Range("B3:AC3").Value2
$snap[4] = $ws.Range("B4:AC4").Value2
$snap[20] = $ws.Range("B20:AC20").Value2
$snap[21] = $ws.Range("B21:AC21").Value2
$snap[23] = $ws.Range("B23:AC23").Value2
$snap[24] = $ws.Range("B24:AC24").Value2
$snap[25] = $ws.Range("B25:AC25").Value2
$snap[26] = $ws.Range("B26:AC26").Value2
$snap[27] = $ws.Range("B27:AC27").Value2
$snap[29] = $ws.Range("B29:AC29").Value2
$snap[43] = $ws.Range("B43:AC43").Value2
$snap[44] = $ws.Range("B44:AC44").Value2
$snap[46] = $ws.Range("B46:AC46").Value2
$snap[47] = $ws.Range("B47:AC47").Value2
$snap[48] = $ws.Range("B48:AC48").Value2
$snap[51] = $ws.Range("B51:AC51").Value2
$snap[52] = $ws.Range("B52:AC52").Value2
$snap[64] = $ws.Range("B64:AC64").Value2
$snap[65] = $ws.Range("B65:AC65").Value2
$snap[70] = $ws.Range("B70:AC70").Value2
$snap[71] = $ws.Range("B71:AC71").Value2
$snap[83] = $ws.Range("B83:AC83").Value2
$snap[84] = $ws.Range("B84:AC84").Value2
$snap[91] = $ws.Range("B91:AC91").Value2
$snap[92] = $ws.Range("B92:AC92").Value2
$snap[116] = $ws.Range("B116:AC116").Value2
$snap[118] = $ws.Range("B118:AC118").Value2
$snap[119] = $ws.Range("B119:AC119").Value2
$snap[120] = $ws.Range("B120:AC120").Value2
$snap[124] = $ws.Range("B124:AC124").Value2
$snap[125] = $ws.Range("B125:AC125").Value2
$snap[132] = $ws.Range("B132:AC132").Value2
$snap[134] = $ws.Range("B134:AC134").Value2
$snap[157] = $ws.Range("B157:AC157").Value2
$snap[158] = $ws.Range("B158:AC158").Value2

# Write the rotated content back: each destination row receives the
# snapshot taken from its source row.
$ws.Range("B3:AC3").Value2 = $snap[4]
$ws.Range("B4:AC4").Value2 = $snap[3]
$ws.Range("B20:AC20").Value2 = $snap[23]
$ws.Range("B21:AC21").Value2 = $snap[24]
$ws.Range("B23:AC23").Value2 = $snap[20]
$ws.Range("B24:AC24").Value2 = $snap[25]
$ws.Range("B25:AC25").Value2 = $snap[21]
$ws.Range("B26:AC26").Value2 = $snap[29]
$ws.Range("B27:AC27").Value2 = $snap[26]
$ws.Range("B29:AC29").Value2 = $snap[27]
$ws.Range("B43:AC43").Value2 = $snap[44]
$ws.Range("B44:AC44").Value2 = $snap[43]
$ws.Range("B46:AC46").Value2 = $snap[47]
$ws.Range("B47:AC47").Value2 = $snap[48]
$ws.Range("B48:AC48").Value2 = $snap[46]
$ws.Range("B51:AC51").Value2 = $snap[52]
$ws.Range("B52:AC52").Value2 = $snap[51]
$ws.Range("B64:AC64").Value2 = $snap[65]
$ws.Range("B65:AC65").Value2 = $snap[64]
$ws.Range("B70:AC70").Value2 = $snap[71]
$ws.Range("B71:AC71").Value2 = $snap[70]
$ws.Range("B83:AC83").Value2 = $snap[84]
$ws.Range("B84:AC84").Value2 = $snap[83]
$ws.Range("B91:AC91").Value2 = $snap[92]
$ws.Range("B92:AC92").Value2 = $snap[91]
$ws.Range("B116:AC116").Value2 = $snap[118]
$ws.Range("B118:AC118").Value2 = $snap[116]
$ws.Range("B119:AC119").Value2 = $snap[120]
$ws.Range("B120:AC120").Value2 = $snap[119]
$ws.Range("B124:AC124").Value2 = $snap[125]
$ws.Range("B125:AC125").Value2 = $snap[124]
$ws.Range("B132:AC132").Value2 = $snap[134]
$ws.Range("B134:AC134").Value2 = $snap[132]
$ws.Range("B157:AC157").Value2 = $snap[158]
$ws.Range("B158:AC158").Value2 = $snap[157]

Write-Output "Row permutation applied."
